# The sheet gains one new data row. A new row is inserted at row 405
# (shifting the existing rows 405-504 down to 406-505) and the new row 405
# is populated with a new price observation for "Pepino ensalada".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 405; this pushes former rows 405..504 down to 406..505
# and picks up the formatting (date number-format on column D) from the row above.
$ws.Rows(405).Insert()

# Populate the newly inserted row 405 with the new record.
$ws.Cells.Item(405, 1).Value  = 5
$ws.Cells.Item(405, 2).Value  = 'Macroferia Regional de Talca'
$ws.Cells.Item(405, 3).Value  = 'Maule'
$ws.Cells.Item(405, 4).Value2 = 44855
$ws.Cells.Item(405, 5).Value  = 7
$ws.Cells.Item(405, 6).Value  = 100112043
$ws.Cells.Item(405, 7).Value  = 'Pepino ensalada'
$ws.Cells.Item(405, 8).Value  = 'Sin especificar'
$ws.Cells.Item(405, 9).Value  = 'Primera'
$ws.Cells.Item(405, 10).Value = 300
$ws.Cells.Item(405, 11).Value = 19000
$ws.Cells.Item(405, 12).Value = 19000
$ws.Cells.Item(405, 13).Value = 19000
$ws.Cells.Item(405, 14).Value = '$/caja 80 unidades'
$ws.Cells.Item(405, 15).Value = 'Región del Maule'
$ws.Cells.Item(405, 16).Value = 238
$ws.Cells.Item(405, 17).Value = 80
$ws.Cells.Item(405, 18).Value = 'Hortaliza'
